# Refatoração da arquitetura do projeto
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Total_de_Uso_por_socio ---
# Lucas Henrique's total usage changes from 1h to 2h.
$ws1 = $wb.Worksheets.Item("Total_de_Uso_por_socio")
$ws1.Range("B5").Value = "2 Horas e 0 Minutos"

# --- Sheet 2: Total_de_Uso_por_Espaço ---
$ws2 = $wb.Worksheets.Item("Total_de_Uso_por_Espaço")
# "Spá" renamed to "Spa"
$ws2.Range("A3").Value = "Spa"
# Insert a new row for "piscina novo" before the "Lago" row (current row 5)
$ws2.Range("A5:B5").Insert()
$ws2.Range("A5").Value = "piscina novo"
$ws2.Range("B5").Value = "1 Horas e 0 Minutos"
# copy style from neighboring data rows so the new row matches formatting
$ws2.Range("A4:B4").Copy()
$ws2.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A5").Value = "piscina novo"
$ws2.Range("B5").Value = "1 Horas e 0 Minutos"

# --- Sheet 3: Total_de_Uso_por_categoria ---
$ws3 = $wb.Worksheets.Item("Total_de_Uso_por_categoria")
# Esportes total usage changes from 10h to 11h
$ws3.Range("B2").Value = "11 Horas e 0 Minutos"
# Relaxamento total usage changes from 3h30 to 4h30
$ws3.Range("B4").Value = "4 Horas e 30 Minutos"
# New row "Lazer" with "0 Horas e 0 Minutos" appended as row 5
$ws3.Range("A5").Value = "Lazer"
$ws3.Range("B5").Value = "0 Horas e 0 Minutos"
